# Generate Report for Handoff
#
# The c8018f2e-b428-4c51-b373-9ec6c5ca8a41.md file has finished translation
# and is now ready to be handed off, for both the zh-cn and de-de locales.
# Update the per-locale status/timestamp rows plus the rollup on Overview.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 is the c8018f2e-... file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("C3").Value = "Ready for handoff"   # de-de status
$wsOverview.Range("D3").Value = "2016-13-12 00:13:14" # Latest Handoff Date

# zh-cn sheet: row 3 is the c8018f2e-... file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"       # Status
$wsZhCn.Range("E3").Value = "2016-03-12 00:13:12"     # Latest Handoff Datetime

# de-de sheet: row 3 is the c8018f2e-... file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"       # Status
$wsDeDe.Range("E3").Value = "2016-03-12 00:13:14"     # Latest Handoff Datetime
